$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 205.88889
$ws.Range("I9").Value = 175.16667
$ws.Range("K9").Value = 175.16667
$ws.Range("M9").Value = -6.166670000000011

$ws.Range("H12").Value = 1229
$ws.Range("I12").Value = 340.2857
$ws.Range("J12").Value = 2265.8333
$ws.Range("K12").Value = 340.2857
$ws.Range("L12").Value = 2265.8333
$ws.Range("M12").Value = -170.2857
$ws.Range("N12").Value = -2605.8333

$ws.Range("H33").Value = 2925.4
$ws.Range("I33").Value = 240.16667
$ws.Range("K33").Value = 240.16667
$ws.Range("M33").Value = -11.16667000000001

$ws.Range("H69").Value = 8423.214
$ws.Range("I69").Value = 8250
$ws.Range("J69").Value = 8452.083000000001
$ws.Range("K69").Value = 24750
$ws.Range("L69").Value = 25356.249
$ws.Range("M69").Value = -23876
$ws.Range("N69").Value = -27104.249

$ws.Range("H72").Value = 8423.214
$ws.Range("I72").Value = 8250
$ws.Range("J72").Value = 8452.083000000001
$ws.Range("K72").Value = 74250
$ws.Range("L72").Value = 76068.747
$ws.Range("M72").Value = -69882
$ws.Range("N72").Value = -84804.747

$ws.Range("H99").Value = 3722
$ws.Range("I99").Value = 3853.875
$ws.Range("J99").Value = 3370.3333
$ws.Range("K99").Value = 11561.625
$ws.Range("L99").Value = 10110.9999
$ws.Range("M99").Value = -10063.625
$ws.Range("N99").Value = -13106.9999

$ws.Range("H106").Value = 64506.5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 64506.5
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 64506.5
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -65768.5

$ws.Range("H138").Value = 3827.48
$ws.Range("I138").Value = 3469.7144
$ws.Range("J138").Value = 4282.8184
$ws.Range("K138").Value = 10409.1432
$ws.Range("L138").Value = 12848.4552
$ws.Range("M138").Value = -5269.143199999999
$ws.Range("N138").Value = -23128.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2242.8
$ws.Range("J3").Value = 1452.5
$ws.Range("L3").Value = 1452.5
$ws.Range("N3").Value = -1682.5

$ws.Range("H4").Value = 1713.5
$ws.Range("I4").Value = 246.5
$ws.Range("J4").Value = 2691.5
$ws.Range("K4").Value = 246.5
$ws.Range("L4").Value = 2691.5
$ws.Range("M4").Value = -130.5
$ws.Range("N4").Value = -2923.5

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H32").Value = 1761.2858
$ws.Range("I32").Value = 1454.0819
$ws.Range("J32").Value = 3843.4443
$ws.Range("K32").Value = 1454.0819
$ws.Range("L32").Value = 3843.4443
$ws.Range("M32").Value = -1167.0819
$ws.Range("N32").Value = -4417.4443

$ws.Range("H45").Value = 71432850
$ws.Range("I45").Value = 90911080
$ws.Range("K45").Value = 90911080
$ws.Range("M45").Value = -90910703

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4999
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2741.5264
$ws.Range("I22").Value = 1144.6364
$ws.Range("J22").Value = 4937.25
$ws.Range("K22").Value = 1144.6364
$ws.Range("L22").Value = 4937.25
$ws.Range("M22").Value = -794.6364000000001
$ws.Range("N22").Value = -5637.25

$ws.Range("H31").Value = 25232.82
$ws.Range("I31").Value = 3157.4707
$ws.Range("K31").Value = 3157.4707
$ws.Range("M31").Value = -2862.4707

$ws.Range("H34").Value = 25232.82
$ws.Range("I34").Value = 3157.4707
$ws.Range("K34").Value = 3157.4707
$ws.Range("M34").Value = -2955.4707

$ws.Range("H58").Value = 4780.143
$ws.Range("I58").Value = 2970.5557
$ws.Range("K58").Value = 2970.5557
$ws.Range("M58").Value = -2767.5557

$ws.Range("H136").Value = 4780.143
$ws.Range("I136").Value = 2970.5557
$ws.Range("K136").Value = 8911.667099999999
$ws.Range("M136").Value = -6361.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 75.5
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 3
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = 174
$ws.Range("N55").Value = -804

$ws.Range("H68").Value = 3518
$ws.Range("J68").Value = 6630
$ws.Range("L68").Value = 19890
$ws.Range("N68").Value = -21512

$ws.Range("H71").Value = 3518
$ws.Range("J71").Value = 6630
$ws.Range("L71").Value = 59670
$ws.Range("N71").Value = -67782

$ws.Range("H131").Value = 6708213
$ws.Range("J131").Value = 7481386
$ws.Range("L131").Value = 22444158
$ws.Range("N131").Value = -22454238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 190.13333
$ws.Range("I2").Value = 136.63637
$ws.Range("J2").Value = 337.25
$ws.Range("K2").Value = 136.63637
$ws.Range("L2").Value = 337.25
$ws.Range("M2").Value = -23.63637
$ws.Range("N2").Value = -563.25

$ws.Range("H97").Value = 1184.4286
$ws.Range("I97").Value = 894.625
$ws.Range("J97").Value = 2111.8
$ws.Range("K97").Value = 894.625
$ws.Range("L97").Value = 2111.8
$ws.Range("M97").Value = -398.625
$ws.Range("N97").Value = -3103.8

$ws.Range("H132").Value = 6952.3794
$ws.Range("I132").Value = 1785.7
$ws.Range("K132").Value = 5357.1
$ws.Range("M132").Value = -2827.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 13500
$ws.Range("J21").Value = 13500
$ws.Range("L21").Value = 13500
$ws.Range("N21").Value = -13848

$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4770

$ws.Range("H93").Value = 43787.312
$ws.Range("I93").Value = 43742.332
$ws.Range("J93").Value = 43814.3
$ws.Range("K93").Value = 43742.332
$ws.Range("L93").Value = 43814.3
$ws.Range("M93").Value = -42494.332
$ws.Range("N93").Value = -46310.3

$ws.Range("H131").Value = 56000
$ws.Range("J131").Value = 56000
$ws.Range("L131").Value = 56000
$ws.Range("N131").Value = -66080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -765
$ws.Range("N21").ClearContents()

$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -710
$ws.Range("N35").ClearContents()

$ws.Range("H58").Value = 6675300
$ws.Range("I58").Value = 13338667
$ws.Range("J58").Value = 11933.333
$ws.Range("K58").Value = 13338667
$ws.Range("L58").Value = 11933.333
$ws.Range("M58").Value = -13338359
$ws.Range("N58").Value = -12549.333

$ws.Range("H61").Value = 22811.4
$ws.Range("I61").Value = 18000
$ws.Range("K61").Value = 18000
$ws.Range("M61").Value = -17708

$ws.Range("H132").Value = 5390.9487
$ws.Range("I132").Value = 4023.3438
$ws.Range("J132").Value = 11642.857
$ws.Range("K132").Value = 12070.0314
$ws.Range("L132").Value = 34928.571
$ws.Range("M132").Value = -9540.0314
$ws.Range("N132").Value = -39988.571

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
